# ZBP_11_obavy_epidemie.xlsx - "aktualizace 1. 6. 2021" update
# Adds a new wave (25. 5. 2021) of survey results as a new trailing column
# on both worksheets, refreshes a handful of previously-published figures
# that were revised upstream, and rewords the two title cells.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")
$ws2 = $wb.Worksheets.Item("pocetR")

# ------------------------------------------------------------------
# Sheet "data" (percentages): new column AE, header = "25. 5. 2021"
# ------------------------------------------------------------------
$ws1.Cells.Item(1, 30).Copy($ws1.Cells.Item(1, 31))   # carry header style (bold + border)
$ws1.Cells.Item(1, 31).Value = "25. 5. 2021"

# A handful of previously published AD (4. 5. 2021) figures were revised
# upstream together with this update:
$ws1.Cells.Item(9, 30).Value = 0.44
$ws1.Cells.Item(10, 30).Value = 0.32
$ws1.Cells.Item(30, 30).Value = 0.4
$ws1.Cells.Item(31, 30).Value = 0.38

# New AE (25. 5. 2021) values, rows 2-76:
$ws1.Cells.Item(2, 31).Value = 0.26
$ws1.Cells.Item(3, 31).Value = 0.43
$ws1.Cells.Item(4, 31).Value = 0.31
$ws1.Cells.Item(5, 31).Value = 0.34
$ws1.Cells.Item(6, 31).Value = 0.45
$ws1.Cells.Item(7, 31).Value = 0.21
$ws1.Cells.Item(8, 31).Value = 0.26
$ws1.Cells.Item(9, 31).Value = 0.47
$ws1.Cells.Item(10, 31).Value = 0.27
$ws1.Cells.Item(11, 31).Value = 0.21
$ws1.Cells.Item(12, 31).Value = 0.38
$ws1.Cells.Item(13, 31).Value = 0.41
$ws1.Cells.Item(14, 31).Value = 0.27
$ws1.Cells.Item(15, 31).Value = 0.43
$ws1.Cells.Item(16, 31).Value = 0.3
$ws1.Cells.Item(17, 31).Value = 0.23
$ws1.Cells.Item(18, 31).Value = 0.41
$ws1.Cells.Item(19, 31).Value = 0.36
$ws1.Cells.Item(20, 31).Value = 0.27
$ws1.Cells.Item(21, 31).Value = 0.44
$ws1.Cells.Item(22, 31).Value = 0.29
$ws1.Cells.Item(23, 31).Value = 0.35
$ws1.Cells.Item(24, 31).Value = 0.41
$ws1.Cells.Item(25, 31).Value = 0.24
$ws1.Cells.Item(26, 31).Value = 0.18
$ws1.Cells.Item(27, 31).Value = 0.45
$ws1.Cells.Item(28, 31).Value = 0.37
$ws1.Cells.Item(29, 31).Value = 0.25
$ws1.Cells.Item(30, 31).Value = 0.41
$ws1.Cells.Item(31, 31).Value = 0.34
$ws1.Cells.Item(32, 31).Value = 0.26
$ws1.Cells.Item(33, 31).Value = 0.45
$ws1.Cells.Item(34, 31).Value = 0.29
$ws1.Cells.Item(35, 31).Value = 0.27
$ws1.Cells.Item(36, 31).Value = 0.43
$ws1.Cells.Item(37, 31).Value = 0.3
$ws1.Cells.Item(38, 31).Value = 0.27
$ws1.Cells.Item(39, 31).Value = 0.44
$ws1.Cells.Item(40, 31).Value = 0.29
$ws1.Cells.Item(41, 31).Value = 0.1
$ws1.Cells.Item(42, 31).Value = 0.44
$ws1.Cells.Item(43, 31).Value = 0.46
$ws1.Cells.Item(44, 31).Value = 0.3
$ws1.Cells.Item(45, 31).Value = 0.46
$ws1.Cells.Item(46, 31).Value = 0.24
$ws1.Cells.Item(47, 31).Value = 0.27
$ws1.Cells.Item(48, 31).Value = 0.44
$ws1.Cells.Item(49, 31).Value = 0.29
$ws1.Cells.Item(50, 31).Value = 0.16
$ws1.Cells.Item(51, 31).Value = 0.53
$ws1.Cells.Item(52, 31).Value = 0.31
$ws1.Cells.Item(53, 31).Value = 0.33
$ws1.Cells.Item(54, 31).Value = 0.43
$ws1.Cells.Item(55, 31).Value = 0.24
$ws1.Cells.Item(56, 31).Value = 0.27
$ws1.Cells.Item(57, 31).Value = 0.44
$ws1.Cells.Item(58, 31).Value = 0.29
$ws1.Cells.Item(59, 31).Value = 0.46
$ws1.Cells.Item(60, 31).Value = 0.39
$ws1.Cells.Item(61, 31).Value = 0.15
$ws1.Cells.Item(62, 31).Value = 0.34
$ws1.Cells.Item(63, 31).Value = 0.42
$ws1.Cells.Item(64, 31).Value = 0.24
$ws1.Cells.Item(65, 31).Value = 0.27
$ws1.Cells.Item(66, 31).Value = 0.4
$ws1.Cells.Item(67, 31).Value = 0.33
$ws1.Cells.Item(68, 31).Value = 0.22
$ws1.Cells.Item(69, 31).Value = 0.51
$ws1.Cells.Item(70, 31).Value = 0.27
$ws1.Cells.Item(71, 31).Value = 0.17
$ws1.Cells.Item(72, 31).Value = 0.51
$ws1.Cells.Item(73, 31).Value = 0.32
$ws1.Cells.Item(74, 31).Value = 0.16
$ws1.Cells.Item(75, 31).Value = 0.36
$ws1.Cells.Item(76, 31).Value = 0.48

# Row 77 repeats the chart title in column A - bump the date there too
$ws1.Cells.Item(77, 1).Value = "Život během pandemie, Obavy z epidemie, % respondentů celkově a ve skupinách, aktualizace 1. 6. 2021"

# ------------------------------------------------------------------
# Sheet "pocetR" (sample sizes): new column AD, header = "25. 5. 2021"
# ------------------------------------------------------------------
$ws2.Cells.Item(1, 29).Copy($ws2.Cells.Item(1, 30))   # carry header style (bold + border)
$ws2.Cells.Item(1, 30).Value = "25. 5. 2021"

# Previously published AC (4. 5. 2021) sample sizes revised upstream:
$ws2.Cells.Item(2, 29).Value = 2029
$ws2.Cells.Item(3, 29).Value = 480
$ws2.Cells.Item(4, 29).Value = 749
$ws2.Cells.Item(5, 29).Value = 800
$ws2.Cells.Item(6, 29).Value = 493
$ws2.Cells.Item(7, 29).Value = 591
$ws2.Cells.Item(9, 29).Value = 980
$ws2.Cells.Item(10, 29).Value = 1049
$ws2.Cells.Item(11, 29).Value = 1058
$ws2.Cells.Item(12, 29).Value = 466
$ws2.Cells.Item(14, 29).Value = 267
$ws2.Cells.Item(19, 29).Value = 271
$ws2.Cells.Item(23, 29).Value = 356
$ws2.Cells.Item(25, 29).Value = 371
$ws2.Cells.Item(26, 29).Value = 443

# New AD (25. 5. 2021) sample-size values, rows 2-26:
$ws2.Cells.Item(2, 30).Value = 1975
$ws2.Cells.Item(3, 30).Value = 471
$ws2.Cells.Item(4, 30).Value = 729
$ws2.Cells.Item(5, 30).Value = 775
$ws2.Cells.Item(6, 30).Value = 443
$ws2.Cells.Item(7, 30).Value = 528
$ws2.Cells.Item(8, 30).Value = 1004
$ws2.Cells.Item(9, 30).Value = 956
$ws2.Cells.Item(10, 30).Value = 1019
$ws2.Cells.Item(11, 30).Value = 1029
$ws2.Cells.Item(12, 30).Value = 456
$ws2.Cells.Item(13, 30).Value = 230
$ws2.Cells.Item(14, 30).Value = 260
$ws2.Cells.Item(15, 30).Value = 53
$ws2.Cells.Item(16, 30).Value = 150
$ws2.Cells.Item(17, 30).Value = 97
$ws2.Cells.Item(18, 30).Value = 17
$ws2.Cells.Item(19, 30).Value = 272
$ws2.Cells.Item(20, 30).Value = 529
$ws2.Cells.Item(21, 30).Value = 239
$ws2.Cells.Item(22, 30).Value = 369
$ws2.Cells.Item(23, 30).Value = 348
$ws2.Cells.Item(24, 30).Value = 232
$ws2.Cells.Item(25, 30).Value = 360
$ws2.Cells.Item(26, 30).Value = 427

# Row 27 is a blank spacer row (every column B..AC is an empty string);
# mirror that into the new AD column by copying a blank neighbour cell.
$ws2.Cells.Item(27, 29).Copy($ws2.Cells.Item(27, 30))

# Row 27 also repeats the chart title in column A - bump the date there too
$ws2.Cells.Item(27, 1).Value = "Život během pandemie, Obavy z epidemie, velikost dotázaného souboru celkově a ve skupinách, aktualizace 1. 6. 2021"

